# Scheduled-runner refresh of market/profit figures across the per-job
# Leve-profit sheets (prices change with each run; other leve metadata
# is untouched).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1092338.9
$ws.Range("I132").Value = 3584.0857
$ws.Range("J132").Value = 4902980.5
$ws.Range("K132").Value = 10752.2571
$ws.Range("L132").Value = 14708941.5
$ws.Range("M132").Value = -8222.257100000001
$ws.Range("N132").Value = -14714001.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 26065.262
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 26065.262
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 26065.262
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -26639.262

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H139").Value = 49112.5
$ws.Range("J139").Value = 49112.5
$ws.Range("L139").Value = 49112.5
$ws.Range("N139").Value = -59392.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 15000
$ws.Range("J61").Value = 15000
$ws.Range("L61").Value = 15000
$ws.Range("N61").Value = -15626

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 47394.95
$ws.Range("J20").Value = 47394.95
$ws.Range("L20").Value = 47394.95
$ws.Range("N20").Value = -47866.95

$ws.Range("H30").Value = 47394.95
$ws.Range("J30").Value = 47394.95
$ws.Range("L30").Value = 47394.95
$ws.Range("N30").Value = -47576.95

$ws.Range("H31").Value = 265365.75
$ws.Range("I31").Value = 79320.38
$ws.Range("J31").Value = 338656.38
$ws.Range("K31").Value = 79320.38
$ws.Range("L31").Value = 338656.38
$ws.Range("M31").Value = -79025.38
$ws.Range("N31").Value = -339246.38

$ws.Range("H34").Value = 265365.75
$ws.Range("I34").Value = 79320.38
$ws.Range("J34").Value = 338656.38
$ws.Range("K34").Value = 79320.38
$ws.Range("L34").Value = 338656.38
$ws.Range("M34").Value = -79118.38
$ws.Range("N34").Value = -339060.38

$ws.Range("H62").Value = 3113.5789
$ws.Range("I62").Value = 2930
$ws.Range("K62").Value = 2930
$ws.Range("M62").Value = -2306

$ws.Range("H65").Value = 3113.5789
$ws.Range("I65").Value = 2930
$ws.Range("K65").Value = 14650
$ws.Range("M65").Value = -11530

$ws.Range("H128").Value = 47394.95
$ws.Range("J128").Value = 47394.95
$ws.Range("L128").Value = 47394.95
$ws.Range("N128").Value = -57354.95

$ws.Range("H132").Value = 38914.85
$ws.Range("I132").Value = 1463.3684
$ws.Range("J132").Value = 127862.125
$ws.Range("K132").Value = 4390.1052
$ws.Range("L132").Value = 383586.375
$ws.Range("M132").Value = -1860.1052
$ws.Range("N132").Value = -388646.375

$ws.Range("H134").Value = 42541.73
$ws.Range("I134").Value = 718.6667
$ws.Range("J134").Value = 55088.65
$ws.Range("K134").Value = 2156.0001
$ws.Range("L134").Value = 165265.95
$ws.Range("M134").Value = 378.9998999999998
$ws.Range("N134").Value = -170335.95

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 2863.6365
$ws.Range("I36").Value = 875
$ws.Range("J36").Value = 4000
$ws.Range("K36").Value = 2625
$ws.Range("L36").Value = 12000
$ws.Range("M36").Value = -2456
$ws.Range("N36").Value = -12338

$ws.Range("H39").Value = 4090
$ws.Range("J39").Value = 4900
$ws.Range("L39").Value = 14700
$ws.Range("N39").Value = -15288

$ws.Range("H107").Value = 1212.9412
$ws.Range("I107").Value = 1462.75
$ws.Range("J107").Value = 990.8889
$ws.Range("K107").Value = 4388.25
$ws.Range("L107").Value = 2972.6667
$ws.Range("M107").Value = -2468.25
$ws.Range("N107").Value = -6812.6667

$ws.Range("H131").Value = 900.14636
$ws.Range("I131").Value = 370.125
$ws.Range("J131").Value = 1028.6364
$ws.Range("K131").Value = 1110.375
$ws.Range("L131").Value = 3085.9092
$ws.Range("M131").Value = 3929.625
$ws.Range("N131").Value = -13165.9092

$ws.Range("H132").Value = 3791.2
$ws.Range("I132").Value = 2044.5714
$ws.Range("J132").Value = 7866.6665
$ws.Range("K132").Value = 18401.1426
$ws.Range("L132").Value = 70799.9985
$ws.Range("M132").Value = -15871.1426
$ws.Range("N132").Value = -75859.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

$ws.Range("H58").Value = 5000
$ws.Range("I58").Value = 5000
$ws.Range("K58").Value = 5000
$ws.Range("M58").Value = -4723

$ws.Range("H97").Value = 2271.6316
$ws.Range("I97").Value = 2459.0908
$ws.Range("K97").Value = 2459.0908
$ws.Range("M97").Value = -1963.0908

$ws.Range("H113").Value = 2235.353
$ws.Range("I113").Value = 1799.75
$ws.Range("J113").Value = 2369.3845
$ws.Range("K113").Value = 1799.75
$ws.Range("L113").Value = 2369.3845
$ws.Range("M113").Value = 370.25
$ws.Range("N113").Value = -6709.3845

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13893477
$ws.Range("I7").Value = 15629286
$ws.Range("J7").Value = 7000
$ws.Range("K7").Value = 15629286
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = -15629174
$ws.Range("N7").Value = -7224

$ws.Range("H61").Value = 1955.8422
$ws.Range("I61").Value = 1969.5518
$ws.Range("J61").Value = 1911.6666
$ws.Range("K61").Value = 1969.5518
$ws.Range("L61").Value = 1911.6666
$ws.Range("M61").Value = -1767.5518
$ws.Range("N61").Value = -2315.6666

$ws.Range("H93").Value = 1071.9395
$ws.Range("I93").Value = 1041.3462
$ws.Range("J93").Value = 1185.5714
$ws.Range("K93").Value = 1041.3462
$ws.Range("L93").Value = 1185.5714
$ws.Range("M93").Value = 206.6538
$ws.Range("N93").Value = -3681.5714

$ws.Range("H100").Value = 1935.9286
$ws.Range("I100").Value = 1775.75
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1775.75
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1234.75
$ws.Range("N100").Value = -3082

$ws.Range("H113").Value = 1955.8422
$ws.Range("I113").Value = 1969.5518
$ws.Range("J113").Value = 1911.6666
$ws.Range("K113").Value = 1969.5518
$ws.Range("L113").Value = 1911.6666
$ws.Range("M113").Value = 200.4482
$ws.Range("N113").Value = -6251.6666

$ws.Range("H126").Value = 13893477
$ws.Range("I126").Value = 15629286
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 46887858
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -46885388
$ws.Range("N126").Value = -25940

$ws.Range("H136").Value = 55433.332
$ws.Range("I136").Value = 38481.035
$ws.Range("K136").Value = 115443.105
$ws.Range("M136").Value = -112893.105

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 9700
$ws.Range("I96").Value = 4000
$ws.Range("J96").Value = 11125
$ws.Range("K96").Value = 4000
$ws.Range("L96").Value = 11125
$ws.Range("M96").Value = -2627
$ws.Range("N96").Value = -13871

$ws.Range("H126").Value = 1129.3182
$ws.Range("I126").Value = 774.8333
$ws.Range("J126").Value = 1554.7
$ws.Range("K126").Value = 2324.4999
$ws.Range("L126").Value = 4664.1
$ws.Range("M126").Value = 145.5001000000002
$ws.Range("N126").Value = -9604.1

$ws.Range("H132").Value = 155806.77
$ws.Range("I132").Value = 112721
$ws.Range("J132").Value = 252749.75
$ws.Range("K132").Value = 338163
$ws.Range("L132").Value = 758249.25
$ws.Range("M132").Value = -335633
$ws.Range("N132").Value = -763309.25
